$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Words" total, which drives the dependent formulas (B4, D9, E9, D11, E11)
$ws.Range("B2").Value = 18631

# Update the active selection to match the saved view state
$ws.Range("D21").Select()
